$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the month typo on earlier diary entries: "Jan" -> "Feb" ---
# These dates were mis-typed when the entries were first authored; the
# diary actually continues into February, so correct the affected rows'
# date column. (The underlying shared date string is reused by more than
# one row, so every row that showed the old "Jan" date is corrected.)
$ws.Range("A29").Value = "5 Feb 2020 (W)"
$ws.Range("A30").Value = "5 Feb 2020 (W)"
$ws.Range("A31").Value = "6 Feb 2020 (Th)"
$ws.Range("A32").Value = "6 Feb 2020 (Th)"
$ws.Range("A33").Value = "12 Feb 2020 (W)"
$ws.Range("A34").Value = "13 Feb 2020 (Th)"

# --- Add the new diary entry for 13 Feb 2020 (row 35) ---
$ws.Range("A35").Value = "13 Feb 2020 (Th)"
$ws.Range("B35").Value = "1445-1700"
$ws.Range("C35").Value = "Harry"
$ws.Range("D35").Formula = "=D34"
$ws.Range("E35").Value = "Went through the final set of slides with the sample questions"
$ws.Range("G35").Value = "I decided not to drink the whole smoothie. Even though it’s advertised as no added sugar, I don’t believe so."

# Row grew to a wrapped multi-line entry, same as the other diary rows.
$ws.Rows.Item(35).RowHeight = 37.3

# --- Move the active selection to the newly-filled cell ---
$ws.Range("E35").Select()
